# Edit script: rewrites the four "How did we run this?" command-line
# paragraphs (English + Hebrew examples) to show the full `java -jar
# ElasticMapReduceRunner.jar ...` invocation and bumps the heb-all run time
# from 39 to 27 minutes; also drops two stray w:hint="cs" paragraph-mark
# hints in the stop-words table row ("167" / "<Hebrew>" cells).

$d = $word.ActiveDocument

$pkgHeader = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body>'
$pkgFooter = '</w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'

function Set-ParagraphXml($para, $bodyXml) {
    $para.Range.InsertXML($pkgHeader + $bodyXml + $pkgFooter)
}

# --- 1) English example paragraph: "0.5 0.2 eng 1 s3://...eng-us-all..." ---
$engBody = '<w:p><w:pPr><w:jc w:val="right"/><w:rPr><w:rFonts w:hint="cs"/><w:u w:val="single"/><w:rtl/></w:rPr></w:pPr><w:proofErr w:type="gramStart"/><w:r><w:rPr><w:rFonts w:cs="Gisha"/><w:sz w:val="16"/><w:szCs w:val="16"/></w:rPr><w:t>java</w:t></w:r><w:proofErr w:type="gramEnd"/><w:r><w:rPr><w:rFonts w:cs="Gisha"/><w:sz w:val="16"/><w:szCs w:val="16"/></w:rPr><w:t xml:space="preserve"> -jar ElasticMapReduceRunner.jar 0.5 0.2 eng 1 s3://datasets.elasticmapreduce/ngrams/books/20090715/eng-us-all/2gram/data</w:t></w:r></w:p>'

# --- 2) heb-all header paragraph: "heb-all [2.4 GB] [252,069,581] - 39 minutes" ---
$hebAllBody = '<w:p><w:pPr><w:jc w:val="right"/><w:rPr><w:sz w:val="32"/><w:szCs w:val="32"/><w:u w:val="single"/><w:rtl/></w:rPr></w:pPr><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:sz w:val="32"/><w:szCs w:val="32"/><w:u w:val="single"/></w:rPr><w:t>heb</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr><w:sz w:val="32"/><w:szCs w:val="32"/><w:u w:val="single"/></w:rPr><w:t>-all [2.4 GB] [252,069,581]</w:t></w:r><w:r><w:rPr><w:sz w:val="32"/><w:szCs w:val="32"/><w:u w:val="single"/></w:rPr><w:t xml:space="preserve"> &#8211; </w:t></w:r><w:r><w:rPr><w:sz w:val="32"/><w:szCs w:val="32"/><w:u w:val="single"/></w:rPr><w:t>27</w:t></w:r><w:r><w:rPr><w:sz w:val="32"/><w:szCs w:val="32"/><w:u w:val="single"/></w:rPr><w:t xml:space="preserve"> minutes</w:t></w:r></w:p>'

# --- 3) Hebrew example paragraph: "0.5 0.2 heb 1 s3://...heb-all..." ---
$hebBody = '<w:p><w:pPr><w:jc w:val="right"/><w:rPr><w:rFonts w:hint="cs"/><w:u w:val="single"/><w:rtl/></w:rPr></w:pPr><w:proofErr w:type="gramStart"/><w:r><w:rPr><w:sz w:val="16"/><w:szCs w:val="16"/></w:rPr><w:t>java</w:t></w:r><w:proofErr w:type="gramEnd"/><w:r><w:rPr><w:sz w:val="16"/><w:szCs w:val="16"/></w:rPr><w:t xml:space="preserve"> -jar ElasticMapReduceRunner.jar 0.5 0.2 </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:sz w:val="16"/><w:szCs w:val="16"/></w:rPr><w:t>heb</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr><w:sz w:val="16"/><w:szCs w:val="16"/></w:rPr><w:t xml:space="preserve"> 1 s3://datasets.elasticmapreduce/ngrams/books/20090715/heb-all/2gram/data</w:t></w:r></w:p>'

$count = $d.Paragraphs.Count
$doneEng = $false
$doneHebAll = $false
$doneHeb = $false

for ($i = 1; $i -le $count; $i++) {
    $para = $d.Paragraphs.Item($i)
    $t = $para.Range.Text

    if (-not $doneEng -and $t -eq "0.5 0.2 eng 1 s3://datasets.elasticmapreduce/ngrams/books/20090715/eng-us-all/2gram/data`r") {
        Set-ParagraphXml $para $engBody
        $doneEng = $true
    }
    elseif (-not $doneHebAll -and $t -eq "heb-all [2.4 GB] [252,069,581] `x{2013} 39 minutes`r") {
        Set-ParagraphXml $para $hebAllBody
        $doneHebAll = $true
    }
    elseif (-not $doneHeb -and $t -eq "0.5 0.2 heb 1 s3://datasets.elasticmapreduce/ngrams/books/20090715/heb-all/2gram/data`r") {
        Set-ParagraphXml $para $hebBody
        $doneHeb = $true
    }
}

Write-Output "eng=$doneEng hebAll=$doneHebAll heb=$doneHeb"

# --- 4) Stop-words table row: drop stray w:hint="cs" on the "167" and
#        Hebrew-holiday-name paragraph marks (no text changes). ---
$cell1Body = '<w:p><w:pPr><w:jc w:val="center"/><w:rPr><w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/><w:sz w:val="18"/><w:szCs w:val="18"/><w:rtl/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial" w:hint="cs"/><w:sz w:val="18"/><w:szCs w:val="18"/><w:rtl/></w:rPr><w:t>167</w:t></w:r></w:p>'
$cell3Body = '<w:p><w:pPr><w:jc w:val="center"/><w:rPr><w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/><w:sz w:val="18"/><w:szCs w:val="18"/><w:rtl/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial" w:hint="cs"/><w:sz w:val="18"/><w:szCs w:val="18"/><w:rtl/></w:rPr><w:t>&#1514;&#1513;&#1506;&#1492; &#1489;&#1488;&#1489;</w:t></w:r></w:p>'

$foundRow = $false
for ($ti = 1; $ti -le $d.Tables.Count; $ti++) {
    $tbl = $d.Tables.Item($ti)
    for ($r = 1; $r -le $tbl.Rows.Count; $r++) {
        $cell1 = $tbl.Cell($r, 1)
        if ($cell1.Range.Text -eq "167`a") {
            $p1 = $cell1.Range.Paragraphs.Item(1)
            Set-ParagraphXml $p1 $cell1Body

            $cell3 = $tbl.Cell($r, 3)
            $p3 = $cell3.Range.Paragraphs.Item(1)
            Set-ParagraphXml $p3 $cell3Body

            $foundRow = $true
        }
    }
}

Write-Output "tableRow=$foundRow"
